# Update the "concise marksheet" figures: correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Number of correct (right) answers changed from 3 to 5
$ws.Range("B11").Value = 5

# Total marks obtained changed from 54 to 90
$ws.Range("B12").Value = 90

# Displayed "obtained/total" string changed from 53/84 to 90/140
$ws.Range("E12").Value = "90/140"
